# "upgrade left table until javakheti": add the 2023 data column (K) to the
# Khulo municipality "Average monthly remuneration" table, continuing the
# yearly series that currently ends at 2022 (column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (year labels) - mirror the formatting of the 2022 header cell
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K3").Value = 2023

# Row 4: total average monthly remuneration
$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("K4").Value = 676.9

# Row 5: women
$ws.Range("J5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("K5").Value = 242.2

# Row 6: men
$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").Value = 703.9

$excel.CutCopyMode = 0
